$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Shift the technique list in column B down by two rows (B4:B29 ->
#    B6:B31), working from the bottom up so we never clobber a value
#    before it has been read. This makes room for two new techniques
#    ("Holden" and "Rizzie Spiral") right after "Spiral5", and folds
#    in the "Thomas Hex" -> "Matthies Hex" rename along the way.
# ------------------------------------------------------------------
for ($r = 29; $r -ge 4; $r--) {
    $val = $ws.Cells.Item($r, 2).Value()
    if ($val -eq "Thomas Hex") {
        $val = "Matthies Hex"
    }
    $ws.Cells.Item($r + 2, 2).Value = $val
}

# ------------------------------------------------------------------
# 2) Insert the two new technique names into the gap created above.
# ------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = "Holden"
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"

# ------------------------------------------------------------------
# 3) Append two brand-new rows (30 and 31) for the extra simulation
#    results, matching the style used by the existing data rows.
# ------------------------------------------------------------------
$ws.Cells.Item(29, 1).Copy($ws.Cells.Item(30, 1))
$ws.Cells.Item(29, 1).Copy($ws.Cells.Item(31, 1))

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Michael-CCHex"

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Michael-SNHex"

$cols = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23)
foreach ($c in $cols) {
    $ws.Cells.Item(30, $c).Value = 1
    $ws.Cells.Item(31, $c).Value = 1
}
